# "scripts aula 02 e 03" - mark attendance for the aula 03 column (E) on the
# attendance sheet ("Plan3"), mirroring what was already entered for aula 02
# (column D) for every student row, then leave the selection where the user
# left off (first cell of the next, still-empty column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan3")

for ($row = 2; $row -le 9; $row++) {
    $aula02 = $ws.Cells.Item($row, 4).Value2   # column D ("aula 02")
    $ws.Cells.Item($row, 5).Value = $aula02    # column E ("aula 03")
}

[void]$ws.Range("E10").Select()
